# Actualización de horarios Línea 141 - 943
# Nuevo timestamp de scraping: 03:51:22 (antes 03:19:42)

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Hoja 1: LP1912
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:51:22"
$ws1.Range("A3").Value = "Total filas: 7"

# Fila 6
$ws1.Cells.Item(6, 1).Value = "03:51:22"
$ws1.Cells.Item(6, 2).Value = "04:02"
$ws1.Cells.Item(6, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(6, 4).Value = 11
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Fila 7
$ws1.Cells.Item(7, 1).Value = "03:51:22"
$ws1.Cells.Item(7, 2).Value = "04:47"
$ws1.Cells.Item(7, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(7, 4).Value = 56
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Fila 8
$ws1.Cells.Item(8, 1).Value = "03:51:22"
$ws1.Cells.Item(8, 2).Value = "04:53"
$ws1.Cells.Item(8, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(8, 4).Value = 62
$ws1.Cells.Item(8, 5).Value = "LP1912"

# Fila 9
$ws1.Cells.Item(9, 1).Value = "03:51:22"
$ws1.Cells.Item(9, 2).Value = "05:11"
$ws1.Cells.Item(9, 3).Value = "17_ROMERO"
$ws1.Cells.Item(9, 4).Value = 80
$ws1.Cells.Item(9, 5).Value = "LP1912"

# Fila 10
$ws1.Cells.Item(10, 1).Value = "03:51:22"
$ws1.Cells.Item(10, 2).Value = "05:22"
$ws1.Cells.Item(10, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(10, 4).Value = 91
$ws1.Cells.Item(10, 5).Value = "LP1912"

# Fila 11 (nueva)
$ws1.Cells.Item(11, 1).Value = "03:51:22"
$ws1.Cells.Item(11, 2).Value = "05:32"
$ws1.Cells.Item(11, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(11, 4).Value = 101
$ws1.Cells.Item(11, 5).Value = "LP1912"

# Fila 12 (nueva)
$ws1.Cells.Item(12, 1).Value = "03:51:22"
$ws1.Cells.Item(12, 2).Value = "05:44"
$ws1.Cells.Item(12, 3).Value = "14_ABASTO"
$ws1.Cells.Item(12, 4).Value = 113
$ws1.Cells.Item(12, 5).Value = "LP1912"

# -----------------------------------------------------------------
# Hoja 2: LP1912-215
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:51:22"

# Fila 6
$ws2.Cells.Item(6, 1).Value = "03:51:22"
$ws2.Cells.Item(6, 2).Value = "04:47"
$ws2.Cells.Item(6, 3).Value = "215_EL PELIGRO"
$ws2.Cells.Item(6, 4).Value = 56
$ws2.Cells.Item(6, 5).Value = "LP1912"

# -----------------------------------------------------------------
# Hoja 3: 6203-6173
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:51:22"
